# daily auto push: 2026-01-18 09:31 UTC
# Insert a new data row at row 649 (shifting existing rows 649:690 down to
# 650:691) and populate it with the new observation for 2026/01/18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 649 downward (existing row 649 and below move to row 650+)
$ws.Rows.Item(649).Insert(-4121)  # -4121 = xlShiftDown

$ws.Cells.Item(649, 1).Value = "'2026/01/18"
$ws.Cells.Item(649, 1).Style = "Normal"
$ws.Cells.Item(649, 2).Value = "日"
$ws.Cells.Item(649, 3).Value = 16
$ws.Cells.Item(649, 4).Value = 23
